# Week 6 update - append new SKU rows (41-46) to the amazon_sales sheet,
# matching the "business_report_week5.xlsx" -> week6 diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append right after the existing last row (40).
# Columns: A=SKU, B=Model, C=(Parent) ASIN, D=(Child) ASIN, P=units_ordered,
# T=ordered_product_sales
$newRows = @(
    @{ Row = 41; A = "FBA79476"; B = "WM-GS1M-BK"; C = "B0DB5VG39T"; D = "B0DB5VG39T"; P = 45; T = 99114.3 },
    @{ Row = 42; A = "FBA79613"; B = "MS1ML";      C = "B0DP2VVRND"; D = "B0DP2VVRND"; P = 40; T = 67423.78 },
    @{ Row = 43; A = "FBA79617"; B = "HDWF1ML";    C = "B0DP32F346"; D = "B0DP32F346"; P = 3;  T = 7855.92 },
    @{ Row = 44; A = "FBA79478"; B = "WM-HA1M-BK"; C = "B0DB5W4TCP"; D = "B0DB5W4TCP"; P = 1;  T = 1397.46 },
    @{ Row = 45; A = "FBA79616"; B = "HD1ML";      C = "B0DP3194QN"; D = "B0DP3194QN"; P = 0;  T = 0 },
    @{ Row = 46; A = "FBA79612"; B = "WM1ML";      C = "B0DP2WC5VW"; D = "B0DP2WC5VW"; P = 0;  T = 0 }
)

# Reuse the currency number format already applied to column T (e.g. T40)
# so the new ordered_product_sales cells keep the same "s=3" style.
$currencyFormat = $ws.Range("T40").NumberFormat

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("P$row").Value = $r.P
    $ws.Range("T$row").Value = $r.T
    $ws.Range("T$row").NumberFormat = $currencyFormat
}

# Update the view: scroll down toward the newly added rows and select B41,
# matching the saved sheetView (topLeftCell="A25", selection B41).
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B41").Select()
